$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster rows (A2:K18) are being re-sorted into a new player order.
# Column A ("No.") stays fixed/sequential; columns B..K (jersey #, name,
# position, height, weight, birth date, nationality, experience, college,
# bbref url) move together as whole rows into their new slots.
#
# Old source row -> new destination row (by player identity):
$mapping = @{
    2  = 2
    3  = 5
    4  = 3
    5  = 4
    6  = 6
    7  = 8
    8  = 7
    9  = 9
    10 = 10
    11 = 13
    12 = 11
    13 = 12
    14 = 14
    15 = 15
    16 = 17
    17 = 16
    18 = 18
}

# Step 1: stash each source row's B:K content into scratch rows far below
# the used range, using Copy so the original cell typing (text vs number,
# shared-string vs inline) is preserved exactly instead of being retyped
# (which would coerce numeric-looking text like "1" into a real number).
$scratchBase = 200
foreach ($srcRow in $mapping.Keys) {
    $scratchRow = $scratchBase + $srcRow
    $ws.Range("B$srcRow`:K$srcRow").Copy($ws.Range("B$scratchRow"))
}

# Step 2: clear out the live player rows so stale cells (e.g. a College
# cell that must end up blank) don't linger once we paste the reordered
# data back in.
foreach ($srcRow in $mapping.Keys) {
    $ws.Range("B$srcRow`:K$srcRow").ClearContents()
}

# Step 3: paste each stashed row into its new destination row.
foreach ($srcRow in $mapping.Keys) {
    $destRow = $mapping[$srcRow]
    $scratchRow = $scratchBase + $srcRow
    $ws.Range("B$scratchRow`:K$scratchRow").Copy($ws.Range("B$destRow"))
}

# Step 4: remove the scratch rows entirely (not just ClearContents) so the
# sheet's used range/dimension shrinks back down to the real data again.
$scratchRows = $mapping.Keys | ForEach-Object { $scratchBase + $_ } | Sort-Object -Descending
foreach ($scratchRow in $scratchRows) {
    $ws.Range("$scratchRow`:$scratchRow").EntireRow.Delete()
}

# Josh Richardson (now row 16) previously had no jersey number on file at
# all; he's been assigned #2 as part of this pass, independent of the
# reshuffle above.
$ws.Range("B16").Value2 = 2
